# Auto-generated edit script: update cryptos list values/percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "51.929.80"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "2.787.55"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'361.93"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'109.75"
$ws.Range("E6").Value = "  -3.63%  "
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = "  -2.30%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "'40.11"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").Value = "'19.48"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "'7.56"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "3.224.66"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "2.784.15"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "51.902.89"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").Value = "'13.14"
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").Value = "'70.38"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'269.61"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'0.161"
$ws.Range("E28").Value = "  +15.22%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.30"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'10.28"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "'51.96"
$ws.Range("E32").Value = "  -3.06%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'18.96"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("E40").Value = "  -4.16%  "
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("D43").Value = "'2.24"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "'119.58"
$ws.Range("E44").Value = "  -6.43%  "
$ws.Range("D45").Value = "'22.00"
$ws.Range("E45").Value = "  -8.19%  "
$ws.Range("D46").Value = "2.084.34"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("D50").Value = "'0.951"
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("D51").Value = "'8.85"
$ws.Range("E51").Value = "  -2.35%  "
